$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 2 entirely (the "hf7hayI1" match). This shifts row 3 (the
# "txqKnEdc" match) up to become the new row 2, matching the diff which
# removes the old row 2 and renumbers the old row 3 as row 2.
$ws.Rows("2").Delete()
